# Aggiunti attributi alle colonne
# Applies the diff: adds a "seriale" column (T) with its datatype/attributes,
# and fills in the "attributi/proprietà" row (now row 6, after inserting a
# blank row 5) plus a new row 7 for additional per-column attributes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a blank row above the old row 5 ("attributi/proprietà" row), ---
# --- shifting it down to row 6 so a new row 7 can hold extra attributes. ---
$ws.Rows("5").Insert()

# --- New column T: "seriale" / VARCHAR(17) ---
$ws.Range("T2").Value = "seriale"
$ws.Range("T2").Font.Bold = $true
$ws.Range("T3").Value = "VARCHAR(17)"

# --- Row 6: attributi/proprietà per column ---
$ws.Range("C6").Value = "PRIMARY_KEY"
$ws.Range("D6").Value = "NOT_NULL"
$ws.Range("E6").Value = "NOT_NULL"
$ws.Range("F6").Value = "NOT_NULL"
$ws.Range("G6").Value = "NOT_NULL"
$ws.Range("H6").Value = "NOT_NULL"
$ws.Range("J6").Value = "NOT_NULL"
$ws.Range("K6").Value = "NOT_NULL"
$ws.Range("N6").Value = "DEFAULT(5)"
$ws.Range("O6").Value = "DEFAULT(False)"
$ws.Range("P6").Value = "DEFAULT(False)"
$ws.Range("Q6").Value = "DEFAULT(False)"
$ws.Range("R6").Value = "DEFAULT(False)"
$ws.Range("S6").Value = "DEFAULT(False)"
$ws.Range("T6").Value = "NOT_NULL"

# --- Row 7: further attributes (defaults / uniqueness) ---
$ws.Range("G7").Value = "DEFAULT(0)"
$ws.Range("H7").Value = "DEFAULT(0)"
$ws.Range("T7").Value = "UNIQUE"

# --- Column width tweaks (COM ColumnWidth snaps to pixel/6-character
# increments, so the inputs below are tuned to land on the closest
# achievable stored width to the target). ---
$ws.Columns("C").ColumnWidth = 12.25
$ws.Columns("Q:S").ColumnWidth = 13.59
$ws.Columns("T").ColumnWidth = 13.1

# --- Final selection, matching the saved workbook state ---
$ws.Range("C9").Select() | Out-Null
